$wb = $excel.ActiveWorkbook

# --- Sheet "Range Status": zero out species counts and clear percentage column ---
$wsRange = $wb.Worksheets.Item("Range Status")
$wsRange.Range("B2").Value = 0
$wsRange.Range("C2").ClearContents()
$wsRange.Range("B3").Value = 0
$wsRange.Range("C3").ClearContents()
$wsRange.Range("B4").Value = 0
$wsRange.Range("C4").ClearContents()
$wsRange.Range("B5").Value = 0
$wsRange.Range("C5").ClearContents()
$wsRange.Range("B6").Value = 0
$wsRange.Range("C6").ClearContents()
$wsRange.Range("B7").Value = 0
$wsRange.Range("C7").ClearContents()

# --- Sheet "Species qualification": Range Analysis count reset to 0 ---
$wsQual = $wb.Worksheets.Item("Species qualification")
$wsQual.Range("B5").Value = 0

# --- Sheet "High Priority break-up": recomputed with only IUCN-based high species ---
$wsBreak = $wb.Worksheets.Item("High Priority break-up")
$wsBreak.Range("A2").Value = "IUCN"
$wsBreak.Range("B2").Value = 8
$wsBreak.Range("C2").Value = 100
$wsBreak.Range("D2").Value = 8
$wsBreak.Range("E2").Value = 100
$wsBreak.Range("A3:E3").ClearContents()

Write-Output "edit applied"
